# "Added user in model 2"
# Adds a second small "model" to the right of the existing permissions
# table: a "Product" header in I1 plus three related labels (Images,
# Reviews, Cart) on row 3 in columns H, I, J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Product"

$ws.Range("H3").Value = "Images"
$ws.Range("I3").Value = "Reviews"
$ws.Range("J3").Value = "Cart"

# Move the active selection to reflect where the user ended up editing.
$ws.Range("J5").Select() | Out-Null
